$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Type"), shifting the old "File Name" column to E.
$ws.Columns.Item(4).Insert()

# Header for the new column.
$ws.Range("D1").Value = "Type"

# Existing voice-line rows (2-31) are all general taunts.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 4).Value = "General"
}

# New rows for recently recorded clips (aggro taunts + a seen/hide sfx pair).
$ws.Range("A32").Value = "`"That's It! I'll kick yer arse!`""
$ws.Range("B32").Value = "X"
$ws.Range("C32").Value = "Me"
$ws.Range("D32").Value = "Aggro"
$ws.Range("E32").Value = "KickYerArse"

$ws.Range("A33").Value = "`"You think you're tough!? I got some scrap for ya!`""
$ws.Range("B33").Value = "X"
$ws.Range("C33").Value = "Me"
$ws.Range("D33").Value = "Aggro"
$ws.Range("E33").Value = "GotScrapForYa"

$ws.Range("A34").Value = "`"Now you see me!`""
$ws.Range("B34").Value = "X"
$ws.Range("C34").Value = "Me"
$ws.Range("D34").Value = "Seen"
$ws.Range("E34").Value = "NowYouSee"

$ws.Range("A35").Value = "`"Now you don't!`""
$ws.Range("B35").Value = "X"
$ws.Range("C35").Value = "Me"
$ws.Range("D35").Value = "Hide"
$ws.Range("E35").Value = "NowYouDont"

# Apply the A1-cell style (used by every row's Quote column) to the new rows.
$ws.Range("A2").Copy()
$ws.Range("A32:A35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Conditional formatting now applies to the full columns instead of a fixed range.
$ws.Range("A1:A1048576").FormatConditions.Delete()
$ws.Range("A1:A1048576").FormatConditions.Add(2, 3, '(INDIRECT("B"&ROW()) = "X")')
$ws.Range("A1:A1048576").FormatConditions.Item(1).Interior.Color = $ws.Range("A1").Interior.Color

$ws.Range("B1:B1048576").FormatConditions.Delete()
$ws.Range("B1:B1048576").FormatConditions.Add(2, 3, '(INDIRECT("B"&ROW())="X")')

$ws.Range("A7").Select()
